$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Cannabis/Marijuana/Weed" (GCDRAN00 / CannabisY/N) and
# "Cocaine powder" (GCDRBN00 / CocaineY/N) feature rows. These are rows
# 42 and 43 in the feature map table; deleting the entire rows shifts
# everything below up by two and keeps the TEXTJOIN formula range (and
# shared-string table) consistent automatically.
$ws.Rows("42:43").Delete() | Out-Null

# Reflect where the author's selection ended up after making the edit.
$ws.Range("G45").Select() | Out-Null
